$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data (and a few reordered coin rows) per the
# "Updated cryptos list" GitHub Actions commit.
#
# The Price/Volume columns are stored as plain text (e.g. "1.00", "26.946.16"),
# not numbers -- force the cells to Text format first so Excel's COM layer
# doesn't silently reinterpret numeric-looking strings (like "1.00") as
# numbers and drop the formatting.
$ws.Range("B2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "26.946.16"
$ws.Range("E2").Value = "  +1.77%  "

# Row 3
$ws.Range("D3").Value = "1.669.75"
$ws.Range("E3").Value = "  +2.78%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "215.79"
$ws.Range("E5").Value = "  +1.04%  "

# Row 6
$ws.Range("D6").Value = "0.546"
$ws.Range("E6").Value = "  +8.92%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").Value = "0.253"
$ws.Range("E8").Value = "  +2.40%  "

# Row 9
$ws.Range("D9").Value = "0.0619"
$ws.Range("E9").Value = "  +1.45%  "

# Row 10
$ws.Range("E10").Value = "  +4.70%  "

# Row 11
$ws.Range("D11").Value = "0.0886"
$ws.Range("E11").Value = "  +3.63%  "

# Row 12
$ws.Range("D12").Value = "1.903.89"
$ws.Range("E12").Value = "  +2.75%  "

# Row 13
$ws.Range("D13").Value = "1.671.27"
$ws.Range("E13").Value = "  +2.10%  "

# Row 14
$ws.Range("D14").Value = "4.07"
$ws.Range("E14").Value = "  +0.52%  "

# Row 15
$ws.Range("D15").Value = "65.78"
$ws.Range("E15").Value = "  +2.65%  "

# Row 16
$ws.Range("E16").Value = "  +1.72%  "

# Row 17
$ws.Range("D17").Value = "26.976.89"
$ws.Range("E17").Value = "  +1.82%  "

# Row 18
$ws.Range("D18").Value = "233.41"
$ws.Range("E18").Value = "  -0.54%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0734"
$ws.Range("E19").Value = "  +1.09%  "

# Row 20
$ws.Range("D20").Value = "7.69"
$ws.Range("E20").Value = "  -0.77%  "

# Row 21
$ws.Range("E21").Value = "  +0.00%  "

# Row 22
$ws.Range("D22").Value = "4.44"
$ws.Range("E22").Value = "  +2.71%  "

# Row 23
$ws.Range("D23").Value = "2.22"
$ws.Range("E23").Value = "  +1.00%  "

# Row 24
$ws.Range("D24").Value = "9.21"
$ws.Range("E24").Value = "  +0.52%  "

# Row 25
$ws.Range("D25").Value = "145.52"
$ws.Range("E25").Value = "  -0.65%  "

# Row 26
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "0.119"
$ws.Range("E26").Value = "  +4.55%  "

# Row 27
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "7.13"
$ws.Range("E27").Value = "  +0.85%  "

# Row 28
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.02%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "15.83"
$ws.Range("E29").Value = "  +1.24%  "

# Row 30
$ws.Range("D30").Value = "0.0497"
$ws.Range("E30").Value = "  +0.29%  "

# Row 31
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +1.26%  "

# Row 32
$ws.Range("E32").Value = "  +1.72%  "

# Row 33
$ws.Range("D33").Value = "1.446.03"
$ws.Range("E33").Value = "  -5.21%  "

# Row 34
$ws.Range("D34").Value = "3.12"
$ws.Range("E34").Value = "  +4.33%  "

# Row 35
$ws.Range("E35").Value = "  +5.36%  "

# Row 36
$ws.Range("E36").Value = "  -0.36%  "

# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "0.892"
$ws.Range("E37").Value = "  +6.74%  "

# Row 38
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.565"
$ws.Range("E38").Value = "  -0.72%  "

# Row 39
$ws.Range("D39").Value = "0.0168"
$ws.Range("E39").Value = "  +1.02%  "

# Row 40
$ws.Range("D40").Value = "6.03"
$ws.Range("E40").Value = "  +2.98%  "

# Row 41
$ws.Range("E41").Value = "  -0.02%  "

# Row 42
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  +9.37%  "

# Row 43
$ws.Range("E43").Value = "  +3.17%  "

# Row 44
$ws.Range("D44").Value = "65.47"
$ws.Range("E44").Value = "  +4.04%  "

# Row 45
$ws.Range("D45").Value = "1.813.65"
$ws.Range("E45").Value = "  +2.80%  "

# Row 46
$ws.Range("D46").Value = "0.778"
$ws.Range("E46").Value = "  +2.27%  "

# Row 47
$ws.Range("D47").Value = "90.49"
$ws.Range("E47").Value = "  +0.74%  "

# Row 48
$ws.Range("D48").Value = "1.52"
$ws.Range("E48").Value = "  +0.90%  "

# Row 49
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0104"
$ws.Range("E49").Value = "  -0.45%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0998"
$ws.Range("E50").Value = "  +3.46%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.0507"
$ws.Range("E51").Value = "  +1.14%  "
